# "update meeting minute 3"
#
# Three semantic changes:
#   1. Center the title paragraph ("Meeting minute 3").
#   2. Move the "_GoBack" bookmark from the end of the "Nhan" attendee
#      paragraph to the very start of the title paragraph (before its run).
#   3. Add an explicit <w:ilvl w:val="0"/> to the <w:numPr> of the
#      paragraphs that currently only carry <w:numId w:val="0"/> (i.e. the
#      "no list" paragraphs scattered through the "What was done" section).

$d = $word.ActiveDocument

# --- 1. Center the "Meeting minute 3" title paragraph ----------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Alignment = 1   # wdAlignParagraphCenter -> <w:jc w:val="center"/>

# --- 2. Relocate the "_GoBack" bookmark -------------------------------------
# A zero-length bookmark range placed exactly at document offset 0 gets
# pushed into the following paragraph by the engine's writer, so we insert a
# throwaway placeholder character first, bookmark right after it (offset 1,
# still inside paragraph 1), then delete the placeholder. Re-adding a
# bookmark with the same name ("_GoBack") automatically removes the
# pre-existing one near "Nhan", so no separate delete step is required there.
$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")

$bookmarkSpot = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$toDelete = $d.Range(0, 1)
$toDelete.Text = ""

# --- 3. Add <w:ilvl w:val="0"/> to the bare numId="0" paragraphs -----------
# These paragraphs only have <w:numPr><w:numId w:val="0"/></w:numPr> today.
# Driving ListLevelNumber through the Word object model also attaches a real
# list (numId becomes nonzero), so immediately calling RemoveNumbers()
# strips the list membership back to numId="0" while leaving the ilvl="0"
# element that Word wrote out in the meantime - exactly the element the
# diff adds.
$noListParagraphs = @(19, 20, 22, 23, 26, 27, 29, 30, 32)
foreach ($idx in $noListParagraphs) {
    $p = $d.Paragraphs($idx)
    $p.Range.ListFormat.ListLevelNumber = 1
    $p.Range.ListFormat.RemoveNumbers()
}
